$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7347
    4  = 5719
    5  = 85
    6  = 179
    10 = 90
    11 = 117
    12 = 208
    13 = 69
    15 = 414
    16 = 53
    17 = 18
    20 = 59
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
